$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$responses = "Q: Handle missing data techniques`nA: dfg d`nQ: Cross-validation methods`nA: g `nQ: Describe complex data analysis project`nA: df `nQ: Stay updated with DS trends`nA: df `nQ: Explain bias-variance tradeoff`nA: df "

$rows = @(14, 15)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "dsfsd"
    $ws.Cells.Item($r, 2).Value = "ib@gmail.com"
    $ws.Cells.Item($r, 3).Value = "Data Scientist"
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 100
    $ws.Cells.Item($r, 6).Value = "Data Visualization, Pytorch, Sql, Machine Learning, Python, Big Data, Tensorflow"
    $ws.Cells.Item($r, 8).Value = $responses
    $ws.Cells.Item($r, 9).Value = "7234f91c-7482-47d9-9870-8cc7ffad8d6b"
}

$ws.Cells.Item(14, 7).Value = "2025-01-26 10:39:28"
$ws.Cells.Item(15, 7).Value = "2025-01-26 10:39:31"
